# Add sample for start_row test
#
# 1. Sheet1 is no longer the active sheet/tab once Sheet3 is appended, and
#    its selection moves to F22.
# 2. A new "Sheet3" is appended after Sheet2 (copy Sheet2 as a template so it
#    picks up the same worksheet-part "look" - x14ac namespace, phoneticPr,
#    pageSetup - then wipe its rows), filled with a 3x8 grid whose cells
#    simply hold their own address ("A1", "B1", ... "C8"), and made the
#    active sheet with selection A10.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("F22").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Copy($null, $ws2)

$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Sheet3"
$ws3.Rows("1:7").Delete()

$columns = @("A", "B", "C")
for ($row = 1; $row -le 8; $row++) {
    foreach ($col in $columns) {
        $address = "$col$row"
        $ws3.Range($address).Value = $address
    }
}

[void]$ws3.Range("A10").Select()
